$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 230
$ws.Range("A3").Value = 233
$ws.Range("A4").Value = 237
$ws.Range("A5").Value = 238.3999999999996
$ws.Range("A6").Value = 226
$ws.Range("A7").Value = 268.5999999999985
$ws.Range("A8").Value = 254.1999999999989
$ws.Range("A9").Value = 293.3999999999996
